# Update of 2025 data and RF changes
# Rows 33 through 74 in column I (RF) change from 17.91404255319149 to 23.09619718309859

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I33:I74").Value = 23.09619718309859
